$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the note text in place first (keeps the shared-string index stable)
$ws.Range("A1").Value = "Note: The time is in seconds"

# Recreate the note row further down (row 15, after the table + average row)
$ws.Range("A15:C15").HorizontalAlignment = -4108
$ws.Range("A15").Value = "Note: The time is in seconds"
$ws.Range("A15:C15").Merge()

# Remove the old note row completely (formatting + contents + merge)
$ws.Range("A1:C1").UnMerge()
$ws.Range("A1:C1").Clear()

# Update the selected cell to match the new target selection
$ws.Range("C19").Select()
